$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.934.75'
$ws.Range('E2').Value = '  -4.45%  '
$ws.Range('D3').Value = '1.738.18'
$ws.Range('E3').Value = '  -4.91%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.50'
$ws.Range('E5').Value = '  -3.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5789'
$ws.Range('E6').Value = '  -3.63%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2731'
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.19'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06615'
$ws.Range('E10').Value = '  -4.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07547'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').Value = '1.742.58'
$ws.Range('E12').Value = '  -4.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.706'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6024'
$ws.Range('E14').Value = '  -4.52%  '
$ws.Range('D15').Value = '1.975.11'
$ws.Range('E15').Value = '  -4.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '74.65'
$ws.Range('E16').Value = '  -3.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008776'
$ws.Range('E17').Value = '  -10.95%  '
$ws.Range('D18').Value = '27.925.74'
$ws.Range('E18').Value = '  -3.67%  '
$ws.Range('E19').Value = '  -4.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '205.57'
$ws.Range('E21').Value = '  -5.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.28'
$ws.Range('E22').Value = '  -2.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.619'
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.28'
$ws.Range('E25').Value = '  -3.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.088'
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('E27').Value = '  -4.44%  '
$ws.Range('E28').Value = '  -2.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.385'
$ws.Range('E29').Value = '  -2.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06145'
$ws.Range('E30').Value = '  -4.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.394'
$ws.Range('E31').Value = '  -3.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.737'
$ws.Range('E32').Value = '  -2.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.728'
$ws.Range('E33').Value = '  -1.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.673'
$ws.Range('E34').Value = '  -2.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.035'
$ws.Range('E35').Value = '  -5.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6415'
$ws.Range('E36').Value = '  -0.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.417'
$ws.Range('E37').Value = '  -5.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.717'
$ws.Range('E38').Value = '  -1.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01668'
$ws.Range('E39').Value = '  -5.23%  '
$ws.Range('D40').Value = '1.126.22'
$ws.Range('E40').Value = '  -1.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.122'
$ws.Range('E41').Value = '  -7.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8749'
$ws.Range('E42').Value = '  -2.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.004'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.71'
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('D45').Value = '1.887.35'
$ws.Range('E45').Value = '  -5.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.29'
$ws.Range('E46').Value = '  -4.64%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.578'
$ws.Range('E47').Value = '  -2.58%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000108'
$ws.Range('E48').Value = '  -4.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.255'
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05377'
$ws.Range('E50').Value = '  -2.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4413'
$ws.Range('E51').Value = '  -2.74%  '
